# ---------------------------------------------------------------------------
# Reproduces the commit "construct target metadata using pydantic model":
#   * appends eleven new metadata worksheets (general, collection, context,
#     licenses, linked_data, meta, provenance, review, sources, spatial,
#     temporal) after the two existing sheets
#   * populates the new "general" sheet with the dataset's Pydantic-modelled
#     metadata (name/topics/title/path/description/... rows, a few external
#     hyperlinks and a publication date)
#   * tweaks the selection/active-tab state on the pre-existing sheets and
#     adds one empty-but-styled cell on the "meta" sheet
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$compSheet = $wb.Worksheets.Item(1)
$metaSheet = $wb.Worksheets.Item(2)

# --- 1. selection bookkeeping on the two original sheets --------------------
$compSheet.Activate()
$compSheet.Range("C8").Select()

$metaSheet.Activate()
$metaSheet.Range("D6").Select()

# a cell that becomes styled (Hyperlink xf) but stays empty
$metaSheet.Range("E4").Style = "Hyperlink"

# --- 2. create the new metadata worksheets, in final tab order -------------
$names = @("general", "collection", "context", "licenses", "linked_data", `
           "meta", "provenance", "review", "sources", "spatial", "temporal")

$newSheets = @{}
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
foreach ($nm in $names) {
    $s = $wb.Worksheets.Add($null, $last)
    $s.Name = $nm
    $newSheets[$nm] = $s
    $last = $s
}

$general = $newSheets["general"]

# --- 3. populate the "general" sheet ----------------------------------------
$general.Columns("A:F").ColumnWidth = 29

# header row
$general.Range("A1").Value = "Propoerty"
$general.Range("B1").Value = "first_or_single"
$general.Range("C1").Value = "second"
$general.Range("D1").Value = "third"
$general.Range("E1").Value = "fourth"

# property rows (column A)
$general.Range("A2").Value  = "name"
$general.Range("A3").Value  = "topics"
$general.Range("A4").Value  = "title"
$general.Range("A5").Value  = "path"
$general.Range("A6").Value  = "description"
$general.Range("A7").Value  = "languages"
$general.Range("A8").Value  = "subject.name"
$general.Range("A9").Value  = "subject.path"
$general.Range("A10").Value = "keywords"
$general.Range("A11").Value = "publicationDate"
$general.Range("A12").Value = "embargoPeriod.start"
$general.Range("A13").Value = "embargoPeriod.end"
$general.Range("A14").Value = "embargoPeriod.isActive"
$general.Range("A2:A14").Style = "Normal"
$general.Range("A2").Font.Bold = $true
$general.Range("A3").Font.Bold = $true
$general.Range("A4").Font.Bold = $true
$general.Range("A5").Font.Bold = $true
$general.Range("A6").Font.Bold = $true
$general.Range("A7").Font.Bold = $true
$general.Range("A8").Font.Bold = $true
$general.Range("A9").Font.Bold = $true
$general.Range("A10").Font.Bold = $true
$general.Range("A11").Font.Bold = $true
$general.Range("A12").Font.Bold = $true
$general.Range("A13").Font.Bold = $true
$general.Range("A14").Font.Bold = $true

# value cells
$general.Range("B2").Value = "Living Lab Measurements"

$general.Range("B3").Value = "Measurement"
$general.Range("C3").Value = "Energy"
$general.Range("D3").Value = "Temperature"

$general.Range("B4").Value = "Living Lab Measurements"

$general.Range("B5").Value = "https://github.com/koubaa-hmc/LLEC_Data/raw/refs/heads/main/dataset_sample_2rows.xlsx"

$general.Range("B6").Value = "The table is a collection of measurements done in a Living Lab"

$general.Range("B7").Value = '"en-GB"'

$general.Range("B8").Value = "energy use"

$general.Range("B9").Value = "http://openenergy-platform.org/ontology/oeo/OEO_00010210"

$general.Range("B10").Value = "http://openenergy-platform.org/ontology/oeo/OEO_00000150"
$general.Range("C10").Value = "http://openenergy-platform.org/ontology/oeo/OEO_00000384"

# publication date (numFmtId 14 - short date)
$general.Range("B11").Value = 45685
$general.Range("B11").NumberFormat = "mm-dd-yy"

# embargoPeriod.isActive
$general.Range("B14").Value = $false

# hyperlinks (added after the date format so it keeps occupying style 7)
$general.Hyperlinks.Add($general.Range("B5"), "https://github.com/koubaa-hmc/LLEC_Data/raw/refs/heads/main/dataset_sample_2rows.xlsx") | Out-Null
$general.Range("B5").Style = "Hyperlink"

$general.Hyperlinks.Add($general.Range("B9"), "http://openenergy-platform.org/ontology/oeo/OEO_00010210") | Out-Null
$general.Range("B9").Style = "Hyperlink"

$general.Hyperlinks.Add($general.Range("B10"), "http://openenergy-platform.org/ontology/oeo/OEO_00000150") | Out-Null
$general.Range("B10").Style = "Hyperlink"

$general.Hyperlinks.Add($general.Range("C10"), "http://openenergy-platform.org/ontology/oeo/OEO_00000384") | Out-Null
$general.Range("C10").Style = "Hyperlink"

# --- 4. final active sheet / selection --------------------------------------
$general.Activate()
$general.Range("B2").Select()
